$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2, 3 and 5 have had their species-record contents rotated:
#   new row2 = old row5
#   new row3 = old row2
#   new row5 = old row3
# Row 4 is untouched. Only columns A, B, D, E, F, G, H, Q, R carry data that
# differs between these three records; capture those before overwriting.

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$row2 = @{}
$row3 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value()
    $row3[$col] = $ws.Range("${col}3").Value()
    $row5[$col] = $ws.Range("${col}5").Value()
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row5[$col]
    $ws.Range("${col}3").Value = $row2[$col]
    $ws.Range("${col}5").Value = $row3[$col]
}
